$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows above the current row 218 (old rows 218:248 shift
# down to become rows 221:251, contents unchanged).
$ws.Rows("218:220").Insert()

# Populate the 3 freshly inserted rows with the new week's data.
$ws.Range("A218:A220").Value = 1
$ws.Range("B218:B220").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C218:C220").Value = "Arica y Parinacota"
$ws.Range("E218:E220").Value = 15
$ws.Range("F218:F220").Value = "Fruta"
$ws.Range("G218:G220").Value = 100106
$ws.Range("H218:H220").Value = "Oleaginosos"
$ws.Range("I218:I220").Value = 100106002
$ws.Range("J218:J220").Value = "Palta"
$ws.Range("K218:K220").Value = "Hass"
$ws.Range("Q218:Q220").Value = "$/bandeja 10 kilos"
$ws.Range("R218:R220").Value = "Perú"
$ws.Range("T218:T220").Value = 10

$ws.Range("D218").Value = 45142
$ws.Range("L218").Value = "Primera"
$ws.Range("M218").Value = 208
$ws.Range("N218").Value = 27000
$ws.Range("O218").Value = 28000
$ws.Range("P218").Value = 27500
$ws.Range("S218").Value = 2750

$ws.Range("D219").Value = 45142
$ws.Range("L219").Value = "Segunda"
$ws.Range("M219").Value = 208
$ws.Range("N219").Value = 25000
$ws.Range("O219").Value = 26000
$ws.Range("P219").Value = 25500
$ws.Range("S219").Value = 2550

$ws.Range("D220").Value = 45142
$ws.Range("L220").Value = "Tercera"
$ws.Range("M220").Value = 104
$ws.Range("N220").Value = 23000
$ws.Range("O220").Value = 24000
$ws.Range("P220").Value = 23500
$ws.Range("S220").Value = 2350
